{"js": "// Bridge Run 6: Add NEAMAP\n// Insert after the \"Run 5: Update maturity\" paragraph:\n//   (blank paragraph)\n//   \"Run 6: Add NEAMAP\"\n//   (blank paragraph)\n// keeping the existing trailing (bookmark-holding) paragraph intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"Run 5: Update maturity\" paragraph.\nlet run5 = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"Run 5: Update maturity\") {\n    run5 = p;\n    break;\n  }\n}\n\nif (!run5) {\n  throw new Error('Could not locate the \"Run 5: Update maturity\" paragraph.');\n}\n\n// Insert the new paragraphs directly after \"Run 5: Update maturity\",\n// in reverse order so each lands immediately below Run 5.\nconst blankAfterRun6 = run5.insertParagraph(\"\", \"After\");\nconst run6 = run5.insertParagraph(\"Run 6: Add NEAMAP\", \"After\");\nconst blankAfterRun5 = run5.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "# Bridge Run 6: Add NEAMAP\n#\n# Target shape (after \"Run 5: Update maturity\"):\n#   Run 5: Update maturity      <- loses the trailing _GoBack bookmark\n#   (blank paragraph)\n#   Run 6: Add NEAMAP\n#   (blank paragraph)\n#   (paragraph that now solely holds the _GoBack bookmark)\n\n$d = $word.ActiveDocument\n\n# Find the \"Run 5: Update maturity\" paragraph and collapse to right after its text.\n$r = $d.Content\n$found = $r.Find.Execute(\"Run 5: Update maturity\")\nif (-not $found) {\n    throw \"Could not find 'Run 5: Update maturity' in the document.\"\n}\n$r.Collapse(0)\n\n# Split the trailing bookmark-only paragraph away from Run 5's own paragraph\n# (this leaves Run 5's paragraph clean and puts the bookmark alone on the\n# document's final paragraph).\n$r.Text = \"`r\"\n$r.Collapse(0)\n\n# Insert a blank paragraph directly after \"Run 5: Update maturity\".\n$r.Text = \"`r\"\n$r.Collapse(0)\n\n# The document's final paragraph is always the bookmark holder. Insert the\n# new \"Run 6: Add NEAMAP\" paragraph (with its own trailing break) right\n# before it.\n$lastIdx = $d.Paragraphs.Count\n$bookmarkPara = $d.Paragraphs.Item($lastIdx)\n$br = $bookmarkPara.Range\n$br.Collapse(1)\n$br.InsertBefore(\"Run 6: Add NEAMAP`r\")\n\n# Finally, split a blank paragraph between \"Run 6: Add NEAMAP\" and the\n# bookmark paragraph.\n$run6Para = $d.Paragraphs.Item($lastIdx)\n$r6 = $run6Para.Range\n$r6.Collapse(0)\n$r6.Text = \"`r\"\n"}
